# "Add files via upload" edit for mosaic-cheatsheet-gf.pptx
#
# The title block's subtitle paragraph currently reads "(ggformula version)"
# (built from three separate runs: "(", "ggformula", " version)" - the
# middle run is flagged err="1", i.e. it was flagged by the spell checker).
# The author simplifies it to a single run reading "ggformula version",
# dropping the surrounding parentheses and collapsing the runs back into
# one, while keeping the original run's character formatting (color/typeface).
#
# Search every slide/shape for the paragraph containing "ggformula" so the
# right shape is found regardless of its shape name/index (the deck reuses
# the name "CustomShape 14" on more than one slide).

$p = $ppt.ActivePresentation

$found = $false
for ($si = 1; $si -le $p.Slides.Count -and -not $found; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count -and -not $found; $i++) {
        $sh = $s.Shapes.Item($i)
        if (-not $sh.HasTextFrame) { continue }

        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -notlike "*ggformula*") { continue }

        $paraCount = $tr.Paragraphs().Count
        for ($pi = 1; $pi -le $paraCount; $pi++) {
            $para = $tr.Paragraphs($pi)
            if ($para.Text -like "*ggformula*") {
                $para.Text = "ggformula version"
                $found = $true
                break
            }
        }
    }
}
